$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "被"
$ws.Range("B7").Value = "passive"
$ws.Range("C7").Value = "虛詞"

$ws.Range("A7:C7").Select()
